# maj fin de journee
# Slide 11 ("Architecture React - Premier Composant") : deux petites retouches
# de texte dans les zones de texte existantes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# ---------------------------------------------------------------------------
# 1) Zone de texte "ZoneTexte 50" (extrait de code) :
#    "const = ComposantTruc = () => {" -> "const ComposantTruc = () => {"
#    (on supprime juste le " = " en trop juste après "const")
# ---------------------------------------------------------------------------
$codeShape = $s.Shapes.Item(20)
$codeRange = $codeShape.TextFrame.TextRange
$codeText  = $codeRange.Text
$eqIndex   = $codeText.IndexOf(" = ")
if ($eqIndex -ge 0) {
    $codeRange.Characters($eqIndex + 1, 3).Text = " "
}

# ---------------------------------------------------------------------------
# 2) Zone de texte "ZoneTexte 1" (puce "La philosophie de react...") :
#    on complète la phrase avec une référence vers le site mdn.
# ---------------------------------------------------------------------------
$philoShape = $s.Shapes.Item(21)
$philoRange = $philoShape.TextFrame.TextRange
$philoText  = $philoRange.Text

$oldTail = " (séparation des concepts) est possible grâce au système import export"
$newTail = $oldTail + " (voir site mdn : référence)"

$tailIndex = $philoText.IndexOf($oldTail)
if ($tailIndex -ge 0) {
    # Réécrit le run existant en une seule fois : le texte ajouté reste donc
    # dans le même run (même mise en forme) que le texte d'origine.
    $philoRange.Characters($tailIndex + 1, $oldTail.Length).Text = $newTail

    # Isole "mdn" dans son propre run, comme dans le fichier d'origine,
    # pour qu'il porte sa propre mise en forme (mot non reconnu -> souligné
    # par le correcteur orthographique).
    $updatedText = $philoRange.Text
    $mdnIndex = $updatedText.IndexOf("mdn")
    if ($mdnIndex -ge 0) {
        $philoRange.Characters($mdnIndex + 1, 3).Text = "mdn"
    }
}
